# NET-FIX, creazione bom sezione obd
# Adds a new "SEZIONE OBD" section to the BOM worksheet (rows 33-61),
# sets a hyperlink on the U1 part-number cell (B44), widens column A,
# and updates the active selection to match the authored workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A33").Value = "SEZIONE OBD"
$ws.Range("B33").Value = "PART NUMBER(NO MOUSER PART…)"
$ws.Range("A35").Value = "connettore OBD"
$ws.Range("D35").Value = 0.9
$ws.Range("A36").Value = "C33-C34"
$ws.Range("B36").Value = "CC0402KRX7R9BB561"
$ws.Range("A37").Value = "R33-R34"
$ws.Range("B37").Value = "ERJ-2RHD1000X"
$ws.Range("A38").Value = "IC7"
$ws.Range("B38").Value = "TJF1051T/3"
$ws.Range("A39").Value = "C31-C32-C27-C23-C24-C25"
$ws.Range("B39").Value = "EMK105BJ105KVHF"
$ws.Range("A40").Value = "IC6"
$ws.Range("B40").Value = "LPC1517JBD48E"
$ws.Range("A41").Value = "C28"
$ws.Range("B41").Value = "04025A100JAT2A"
$ws.Range("A42").Value = "Y2"
$ws.Range("B42").Value = "ECS-120-10-36B-CWY-TR"
$ws.Range("A43").Value = "C30"
$ws.Range("B43").Value = "04025A100JAT2A"
$ws.Range("A44").Value = "U1"
$ws.Range("B44").Value = "AP7313-33SRG-7"
$ws.Range("A45").Value = "R35"
$ws.Range("B45").Value = "CRCW040247K0FKED"
$ws.Range("A46").Value = "R36-R37"
$ws.Range("B46").Value = "TNPW040210K0DEED"
$ws.Range("A47").Value = "C35-C36-C26"
$ws.Range("B47").Value = "C1005X5R1H104K050BB"
$ws.Range("A48").Value = "P1"
$ws.Range("B48").Value = "TLV76050DBZT"
$ws.Range("A49").Value = "R38-R39"
$ws.Range("B49").Value = "CRCW06030000Z0EC"
$ws.Range("A50").Value = "U2<A:B>"
$ws.Range("B50").Value = "LM393D"
$ws.Range("A51").Value = "R40-R41"
$ws.Range("B51").Value = "RC0402FR-13470RL"
$ws.Range("A52").Value = "DS1"
$ws.Range("B52").Value = "HSMS-C190"
$ws.Range("A53").Value = "DS2"
$ws.Range("B53").Value = "HSMY-C190"
$ws.Range("A54").Value = "R32"
$ws.Range("B54").Value = "CRCW0402510RFKEDC"
$ws.Range("A55").Value = "U3"
$ws.Range("B55").Value = "MCP2021-500E/SN"
$ws.Range("A56").Value = "C29"
$ws.Range("B56").Value = "GRM155R60J106ME05J"
$ws.Range("A57").Value = "R26-R27"
$ws.Range("B57").Value = "ERJ-2RKF3301X"
$ws.Range("A58").Value = "R28-R29-R31-R30"
$ws.Range("B58").Value = "TNPW040210K0DEED"
$ws.Range("A59").Value = "D4"
$ws.Range("B59").Value = "BAT46W-7-F"
$ws.Range("A60").Value = "Q5"
$ws.Range("B60").Value = "BSS84LT1G"
$ws.Range("A61").Value = "Q6-Q7"
$ws.Range("B61").Value = "BSS138"

# Hyperlink on U1's part number (AP7313-33SRG-7)
$ws.Hyperlinks.Add($ws.Range("B44"), "https://octopart.com/ap7313-33srg-7-diodes+inc.-17931362?r=sp", "", "", "https://octopart.com/ap7313-33srg-7-diodes+inc.-17931362?r=sp")

# Widen column A to fit the new, longer component labels
$ws.Columns.Item(1).ColumnWidth = 27.6

# Restore the selection/viewport state recorded in the authored workbook
$ws.Range("B63").Select()
